# Auto-generated script applying the Cactuar_Profits (Sheets workbook) diff.
# Each sheet corresponds to a Leve-crafting profit table; only literal <v> numeric
# cells changed (no formulas in this workbook), so we write plain Range.Value
# assignments, and ClearContents() for cells whose <c> element was removed outright.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 763.0833
$ws.Range("I33").Value = 763.0833
$ws.Range("K33").Value = 763.0833
$ws.Range("M33").Value = -534.0833
$ws.Range("H40").Value = 15647760
$ws.Range("J40").Value = 33365370
$ws.Range("L40").Value = 33365370
$ws.Range("N40").Value = -33365720
$ws.Range("H43").Value = 1713556.2
$ws.Range("J43").Value = 4627.2856
$ws.Range("L43").Value = 4627.2856
$ws.Range("N43").Value = -4765.2856
$ws.Range("H82").Value = 4745.9165
$ws.Range("I82").Value = 993
$ws.Range("K82").Value = 2979
$ws.Range("M82").Value = -2573
$ws.Range("H85").Value = 4745.9165
$ws.Range("I85").Value = 993
$ws.Range("K85").Value = 2979
$ws.Range("M85").Value = -1575
$ws.Range("H106").Value = 20835518
$ws.Range("J106").Value = 4653
$ws.Range("L106").Value = 4653
$ws.Range("N106").Value = -5915
$ws.Range("H116").Value = 69280000
$ws.Range("I116").Value = 60842596
$ws.Range("J116").Value = 83342330
$ws.Range("K116").Value = 60842596
$ws.Range("L116").Value = 83342330
$ws.Range("M116").Value = -60839154
$ws.Range("N116").Value = -83349214
$ws.Range("H129").Value = 1829.1765
$ws.Range("I129").Value = 985.8570999999999
$ws.Range("K129").Value = 2957.5713
$ws.Range("M129").Value = 2042.4287
$ws.Range("H137").Value = 9282443
$ws.Range("J137").Value = 22228860
$ws.Range("L137").Value = 66686580
$ws.Range("N137").Value = -66691680
$ws.Range("H138").Value = 2665.625
$ws.Range("I138").Value = 2554.1667
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 7662.500100000001
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = -2522.500100000001
$ws.Range("N138").Value = -19280

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15334.178
$ws.Range("I32").Value = 15669.823
$ws.Range("K32").Value = 15669.823
$ws.Range("M32").Value = -15382.823
$ws.Range("H132").Value = 17464.236
$ws.Range("I132").Value = 22197.814
$ws.Range("J132").Value = 5845.4546
$ws.Range("K132").Value = 66593.442
$ws.Range("L132").Value = 17536.3638
$ws.Range("M132").Value = -64063.442
$ws.Range("N132").Value = -22596.3638

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1766.0834
$ws.Range("I86").Value = 2032.5
$ws.Range("J86").Value = 1499.6666
$ws.Range("K86").Value = 2032.5
$ws.Range("L86").Value = 1499.6666
$ws.Range("M86").Value = -909.5
$ws.Range("N86").Value = -3745.6666
$ws.Range("H89").Value = 1766.0834
$ws.Range("I89").Value = 2032.5
$ws.Range("J89").Value = 1499.6666
$ws.Range("K89").Value = 10162.5
$ws.Range("L89").Value = 7498.333000000001
$ws.Range("M89").Value = -4546.5
$ws.Range("N89").Value = -18730.333
$ws.Range("H99").Value = 2605113.8
$ws.Range("J99").Value = 1142.25
$ws.Range("L99").Value = 1142.25
$ws.Range("N99").Value = -4138.25
$ws.Range("H105").Value = 90910840
$ws.Range("I105").Value = 111112590
$ws.Range("K105").Value = 111112590
$ws.Range("M105").Value = -111110843

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 921.5
$ws.Range("I22").Value = 637.75
$ws.Range("J22").Value = 1299.8334
$ws.Range("K22").Value = 637.75
$ws.Range("L22").Value = 1299.8334
$ws.Range("M22").Value = -287.75
$ws.Range("N22").Value = -1999.8334
$ws.Range("H31").Value = 5259.4106
$ws.Range("I31").Value = 1919.8
$ws.Range("K31").Value = 1919.8
$ws.Range("M31").Value = -1624.8
$ws.Range("H34").Value = 5259.4106
$ws.Range("I34").Value = 1919.8
$ws.Range("K34").Value = 1919.8
$ws.Range("M34").Value = -1717.8
$ws.Range("H132").Value = 11914078
$ws.Range("I132").Value = 14503599
$ws.Range("K132").Value = 43510797
$ws.Range("M132").Value = -43508267
$ws.Range("H134").Value = 2030.2106
$ws.Range("I134").Value = 1739.6471
$ws.Range("K134").Value = 5218.9413
$ws.Range("M134").Value = -2683.9413

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 91201290
$ws.Range("J11").Value = 273336670
$ws.Range("L11").Value = 820010010
$ws.Range("N11").Value = -820010290
$ws.Range("H37").Value = 114990
$ws.Range("J37").Value = 114990
$ws.Range("L37").Value = 344970
$ws.Range("N37").Value = -345194
$ws.Range("H75").Value = 6081.6
$ws.Range("J75").Value = 7014.5
$ws.Range("L75").Value = 21043.5
$ws.Range("N75").Value = -23039.5
$ws.Range("H78").Value = 6081.6
$ws.Range("J78").Value = 7014.5
$ws.Range("L78").Value = 63130.5
$ws.Range("N78").Value = -73114.5
$ws.Range("H113").Value = 737.2222
$ws.Range("J113").Value = 757.8182
$ws.Range("L113").Value = 2273.4546
$ws.Range("N113").Value = -6613.4546
$ws.Range("H116").Value = 10479.2
$ws.Range("I116").Value = 5000
$ws.Range("J116").Value = 11849
$ws.Range("K116").Value = 15000
$ws.Range("L116").Value = 35547
$ws.Range("M116").Value = -11558
$ws.Range("N116").Value = -42431
$ws.Range("H140").Value = 5492.409
$ws.Range("I140").Value = 4118.6
$ws.Range("J140").Value = 6637.25
$ws.Range("K140").Value = 12355.8
$ws.Range("L140").Value = 19911.75
$ws.Range("M140").Value = -7175.800000000001
$ws.Range("N140").Value = -30271.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 346.05
$ws.Range("I2").Value = 135.4
$ws.Range("K2").Value = 135.4
$ws.Range("M2").Value = -22.40000000000001
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H70").Value = 2511242
$ws.Range("J70").Value = 5007.4
$ws.Range("L70").Value = 5007.4
$ws.Range("N70").Value = -5547.4
$ws.Range("H73").Value = 2511242
$ws.Range("J73").Value = 5007.4
$ws.Range("L73").Value = 5007.4
$ws.Range("N73").Value = -6879.4
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H107").Value = 7937341
$ws.Range("J107").Value = 999.3333
$ws.Range("L107").Value = 999.3333
$ws.Range("N107").Value = -4839.3333
$ws.Range("H132").Value = 117529.945
$ws.Range("I132").Value = 185963.1
$ws.Range("K132").Value = 557889.3
$ws.Range("M132").Value = -555359.3

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 42859024
$ws.Range("I16").Value = 60001396
$ws.Range("J16").Value = 3087.25
$ws.Range("K16").Value = 60001396
$ws.Range("L16").Value = 3087.25
$ws.Range("M16").Value = -60001226
$ws.Range("N16").Value = -3427.25
$ws.Range("H55").Value = 368.5625
$ws.Range("I55").Value = 136.2
$ws.Range("J55").Value = 474.18182
$ws.Range("K55").Value = 136.2
$ws.Range("L55").Value = 474.18182
$ws.Range("M55").Value = 36.80000000000001
$ws.Range("N55").Value = -820.18182
$ws.Range("H101").Value = 58998.75
$ws.Range("J101").Value = 58998.75
$ws.Range("L101").Value = 58998.75
$ws.Range("N101").Value = -65488.75
$ws.Range("H132").Value = 4546.0293
$ws.Range("J132").Value = 6489.75
$ws.Range("L132").Value = 19469.25
$ws.Range("N132").Value = -24529.25
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 4898.6665
$ws.Range("I136").Value = 3484.2856
$ws.Range("K136").Value = 10452.8568
$ws.Range("M136").Value = -7902.856800000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4638857
$ws.Range("J81").Value = 3482388.2
$ws.Range("L81").Value = 6964776.4
$ws.Range("N81").Value = -6966898.4
$ws.Range("H84").Value = 4638857
$ws.Range("J84").Value = 3482388.2
$ws.Range("L84").Value = 34823882
$ws.Range("N84").Value = -34834490
$ws.Range("H126").Value = 2334.85
$ws.Range("I126").Value = 2682.25
$ws.Range("K126").Value = 8046.75
$ws.Range("M126").Value = -5576.75
